$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.157.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.45%  '

$ws.Range("D3").Value = "'1.656.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.17%  '

$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").Value = "'217.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.91%  '

$ws.Range("D6").Value = "'0.5149"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.48%  '

$ws.Range("D7").Value = "'1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("E8").Value = '  -3.77%  '

$ws.Range("D9").Value = "'0.06404"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.08%  '

$ws.Range("D10").Value = "'19.83"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.07778"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.78%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'4.305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.55%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'1.654.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.84%  '

$ws.Range("D14").Value = "'1.884.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("D15").Value = "'0.5528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.10%  '

$ws.Range("D16").Value = "'0.0₅8033"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.68%  '

$ws.Range("D17").Value = "'64.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.25%  '

$ws.Range("D18").Value = "'26.194.47"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("D20").Value = "'210.34"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'4.388"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.87%  '

$ws.Range("D22").Value = "'10.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.76%  '

$ws.Range("D23").Value = "'5.878"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.73%  '

$ws.Range("E24").Value = '  +0.29%  '

$ws.Range("D25").Value = "'144.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.23%  '

$ws.Range("D26").Value = "'1.759"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.54%  '

$ws.Range("D27").Value = "'0.1160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.60%  '

$ws.Range("D28").Value = "'6.955"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.46%  '

$ws.Range("D29").Value = "'15.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.25%  '

$ws.Range("D30").Value = "'0.05251"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.75%  '

$ws.Range("D31").Value = "'1.256"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.77%  '

$ws.Range("D32").Value = "'3.366"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.51%  '

$ws.Range("D33").Value = "'3.206"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.42%  '

$ws.Range("D34").Value = "'1.563"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.80%  '

$ws.Range("D36").Value = "'2.366"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.08%  '

$ws.Range("D37").Value = "'0.9236"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.64%  '

$ws.Range("D38").Value = "'0.5722"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.94%  '

$ws.Range("D39").Value = "'1.162.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.16%  '

$ws.Range("D40").Value = "'0.01587"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.77%  '

$ws.Range("E41").Value = '  +0.34%  '

$ws.Range("D42").Value = "'0.8387"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("D43").Value = "'5.661"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'99.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.90%  '

$ws.Range("D45").Value = "'1.794.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.23%  '

$ws.Range("D46").Value = "'0.0₈113"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.22%  '

$ws.Range("D47").Value = "'0.4501"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("D48").Value = "'56.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.41%  '

$ws.Range("E49").Value = '  +0.69%  '

$ws.Range("D50").Value = "'7.901"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.13%  '

$ws.Range("D51").Value = "'0.05095"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.80%  '
